$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section header row (clone formatting of row 2) ---
$ws.Range("A2:E2").Copy($ws.Range("A25:E25"))
$ws.Range("A25").Value = "TC # 01.02 - Test"

# --- Backfill the "TC #" label into column A for the existing
#     "TC # 01.01 - Shop by search and category" block (rows 3-23) ---
for ($r = 3; $r -le 23; $r++) {
    $ws.Cells.Item($r, 1).Value = "TC # 01.01"
}

# --- New block of steps (clone formatting+values of rows 3-23) ---
$ws.Range("A3:E23").Copy($ws.Range("A26:E46"))
for ($r = 26; $r -le 46; $r++) {
    $ws.Cells.Item($r, 1).Value = "TC # 01.02"
}

# --- New footer row (clone formatting+values of row 24) ---
$ws.Range("A24:E24").Copy($ws.Range("A47:E47"))

# --- View / layout tweaks ---
$ws.Columns(1).ColumnWidth = 9
[void]$ws.Range("E50").Select()
